# Pooh Points: normal 20260201
# Refresh live game stats on the Players sheet and roll the totals on OwnerTotals.
$wb = $excel.ActiveWorkbook

# --- Players sheet: live box-score refresh ---
$ws = $wb.Worksheets.Item("Players")

# Status text got shorter ("18:55 - 2nd Half" vs a full scheduled-time string) -> narrow column G
$ws.Columns.Item(7).ColumnWidth = 17.142857142857142

# Row 13
$ws.Range("G13").Value = '18:55 - 2nd Half'
$ws.Range("H13").Value = 3
$ws.Range("I13").Value = 11
$ws.Range("K13").Value = 1
$ws.Range("N13").Value = 4
$ws.Range("P13").Value = 25
$ws.Range("Q13").Value = 4
$ws.Range("R13").Value = 10
$ws.Range("S13").Value = 3
$ws.Range("T13").Value = 4

# Row 23
$ws.Range("G23").Value = '18:55 - 2nd Half'
$ws.Range("H23").Value = 17
$ws.Range("I23").Value = 19
$ws.Range("J23").Value = 3
$ws.Range("L23").Value = 1
$ws.Range("N23").Value = 1
$ws.Range("P23").Value = 28
$ws.Range("Q23").Value = 8
$ws.Range("R23").Value = 15
$ws.Range("T23").Value = 6

# Row 24
$ws.Range("G24").Value = '18:55 - 2nd Half'
$ws.Range("H24").Value = 4
$ws.Range("I24").Value = 2
$ws.Range("J24").Value = 2
$ws.Range("P24").Value = 9
$ws.Range("Q24").Value = 1
$ws.Range("R24").Value = 1

# Row 32
$ws.Range("G32").Value = '18:55 - 2nd Half'
$ws.Range("H32").Value = 18
$ws.Range("I32").Value = 19
$ws.Range("J32").Value = 3
$ws.Range("K32").Value = 2
$ws.Range("N32").Value = 2
$ws.Range("P32").Value = 24
$ws.Range("Q32").Value = 7
$ws.Range("R32").Value = 11
$ws.Range("T32").Value = 6
$ws.Range("U32").Value = 2
$ws.Range("V32").Value = 2

# Row 33
$ws.Range("G33").Value = '18:55 - 2nd Half'

# Row 43
$ws.Range("G43").Value = '18:55 - 2nd Half'
$ws.Range("H43").Value = 31
$ws.Range("I43").Value = 21
$ws.Range("J43").Value = 6
$ws.Range("K43").Value = 4
$ws.Range("L43").Value = 2
$ws.Range("M43").Value = 2
$ws.Range("O43").Value = 2
$ws.Range("P43").Value = 27
$ws.Range("Q43").Value = 9
$ws.Range("R43").Value = 12
$ws.Range("T43").Value = 1
$ws.Range("U43").Value = 3
$ws.Range("V43").Value = 4

# Row 44
$ws.Range("G44").Value = '18:55 - 2nd Half'
$ws.Range("H44").Value = 18
$ws.Range("I44").Value = 10
$ws.Range("K44").Value = 7
$ws.Range("L44").Value = 6
$ws.Range("O44").Value = 2
$ws.Range("P44").Value = 27
$ws.Range("Q44").Value = 5
$ws.Range("R44").Value = 11
$ws.Range("T44").Value = 1

# Row 45
$ws.Range("G45").Value = '18:55 - 2nd Half'
$ws.Range("H45").Value = 1
$ws.Range("I45").Value = 3
$ws.Range("J45").Value = 3
$ws.Range("P45").Value = 17
$ws.Range("Q45").Value = 1
$ws.Range("R45").Value = 4
$ws.Range("S45").Value = 1
$ws.Range("T45").Value = 3

# Row 55
$ws.Range("G55").Value = '18:55 - 2nd Half'
$ws.Range("H55").Value = 4
$ws.Range("I55").Value = 6
$ws.Range("J55").Value = 3
$ws.Range("O55").Value = 3
$ws.Range("P55").Value = 17
$ws.Range("Q55").Value = 3
$ws.Range("R55").Value = 11
$ws.Range("T55").Value = 3

# Row 56
$ws.Range("G56").Value = '18:55 - 2nd Half'
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 3
$ws.Range("K56").Value = 2
$ws.Range("L56").Value = 1
$ws.Range("N56").Value = 3
$ws.Range("O56").Value = 4
$ws.Range("P56").Value = 17
$ws.Range("R56").Value = 6
$ws.Range("T56").Value = 3

# Row 65
$ws.Range("G65").Value = '18:55 - 2nd Half'
$ws.Range("H65").Value = 17
$ws.Range("I65").Value = 7
$ws.Range("J65").Value = 13
$ws.Range("M65").Value = 2
$ws.Range("O65").Value = 1
$ws.Range("P65").Value = 23
$ws.Range("Q65").Value = 2
$ws.Range("R65").Value = 9
$ws.Range("U65").Value = 3
$ws.Range("V65").Value = 4

# Row 66
$ws.Range("G66").Value = '18:55 - 2nd Half'
$ws.Range("H66").Value = 6
$ws.Range("J66").Value = 4
$ws.Range("O66").Value = 1
$ws.Range("P66").Value = 20

# Row 78
$ws.Range("G78").Value = '18:55 - 2nd Half'
$ws.Range("H78").Value = 12
$ws.Range("I78").Value = 8
$ws.Range("J78").Value = 4
$ws.Range("K78").Value = 3
$ws.Range("O78").Value = 1
$ws.Range("P78").Value = 22
$ws.Range("Q78").Value = 3
$ws.Range("R78").Value = 5

# Row 146
$ws.Range("G146").Value = '18:55 - 2nd Half'
$ws.Range("H146").Value = 11
$ws.Range("I146").Value = 6
$ws.Range("J146").Value = 6
$ws.Range("M146").Value = 4
$ws.Range("N146").Value = 2
$ws.Range("O146").Value = 3
$ws.Range("P146").Value = 21
$ws.Range("Q146").Value = 2
$ws.Range("R146").Value = 4
$ws.Range("V146").Value = 3

# Row 147
$ws.Range("G147").Value = '18:55 - 2nd Half'
$ws.Range("H147").Value = 8
$ws.Range("I147").Value = 5
$ws.Range("J147").Value = 3
$ws.Range("K147").Value = 3
$ws.Range("O147").Value = 3
$ws.Range("P147").Value = 20
$ws.Range("Q147").Value = 2
$ws.Range("R147").Value = 5
$ws.Range("S147").Value = 1
$ws.Range("T147").Value = 2

# Row 148
$ws.Range("G148").Value = '18:55 - 2nd Half'
$ws.Range("O148").Value = 2
$ws.Range("P148").Value = 10

# Row 149
$ws.Range("D149").Value = 'London Jemison'
$ws.Range("G149").Value = '18:55 - 2nd Half'
$ws.Range("H149").Value = 2
$ws.Range("I149").Value = 3
$ws.Range("J149").Value = 0
$ws.Range("M149").Value = 1
$ws.Range("P149").Value = 8
$ws.Range("Q149").Value = 1
$ws.Range("R149").Value = 3
$ws.Range("S149").Value = 1
$ws.Range("T149").Value = 2

# Row 150
$ws.Range("D150").Value = 'Noah Williamson'
$ws.Range("G150").Value = '18:55 - 2nd Half'
$ws.Range("H150").Value = 1
$ws.Range("J150").Value = 1
$ws.Range("P150").Value = 5
$ws.Range("R150").Value = 0

# --- OwnerTotals sheet: updated point totals (ranking reshuffled) ---
$ws2 = $wb.Worksheets.Item("OwnerTotals")

# Row 2
$ws2.Range("B2").Value = 96

# Row 3
$ws2.Range("B3").Value = 84

# Row 4
$ws2.Range("B4").Value = 74

# Row 5
$ws2.Range("A5").Value = 'G-Flop'
$ws2.Range("B5").Value = 59
$ws2.Range("C5").Value = 5

# Row 6
$ws2.Range("A6").Value = 'Hilton Heads'
$ws2.Range("B6").Value = 56
$ws2.Range("C6").Value = 4

# Row 7
$ws2.Range("B7").Value = 52
